$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Latest period (release date)" value for the
# "Employment share by occupation" row (C3), which had mistakenly
# been left the same as the row above it.
$ws.Range("C3").Value = "Jan 2021 - Dec 2021 (12/04/22)"

# Update the selected cell to reflect where the edit was made.
$ws.Range("C4").Select()
